$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.850754976272583
$ws.Range("B1").Value = 1.566083669662476
$ws.Range("C1").Value = 7.160388946533203
$ws.Range("D1").Value = 2.683520317077637
$ws.Range("E1").Value = 1.574187874794006
